$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 4 (Fecha, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Precio $/Kg)

$ws.Range("D2").Value = 44414
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 31000
$ws.Range("L2").Value = 32000
$ws.Range("M2").Value = 31500
$ws.Range("P2").Value = 1260

$ws.Range("D4").Value = 44827
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 31000
$ws.Range("M4").Value = 30500
$ws.Range("P4").Value = 1220
